# export timeseries functions to excel (60/1114)
#
# Adds a new "TimeSeries" unit-test sheet (between "Date" and "Utilities"),
# lists it in the "UnitTests" summary sheet, flips two previously-failing
# "Date" sheet assertions to PASS, and leaves "Utilities" as the
# last / active tab.

$wb = $excel.ActiveWorkbook

$dateSheet = $wb.Worksheets.Item("Date")
$utilSheet = $wb.Worksheets.Item("Utilities")
$unitTests = $wb.Worksheets.Item("UnitTests")

# ---------------------------------------------------------------------
# 1. Two Date-sheet regression tests that used to fail now pass.
#    Row 30 -> qlECBKnownDates ; Row 39 -> qlECBIsECBdate
# ---------------------------------------------------------------------
$dateSheet.Cells.Item(30, 3).Value = 25569
$dateSheet.Cells.Item(30, 5).Value = 25569

$dateSheet.Cells.Item(39, 3).Value = $true
$dateSheet.Cells.Item(39, 5).Value = $true

# ---------------------------------------------------------------------
# 2. Insert the new "TimeSeries" worksheet right before "Utilities".
# ---------------------------------------------------------------------
$ts = $wb.Worksheets.Add($utilSheet)
$ts.Name = "TimeSeries"

# Header rows (same layout as every other unit-test sheet).
$ts.Range("A1").Value = "Function"
$ts.Range("B1").Value = "Expected"
$ts.Range("C1").Value = "Actual"
$ts.Range("D1").Value = "PASS /"
$ts.Range("E1").Value = "Function"

$ts.Range("A2").Value = "Name"
$ts.Range("B2").Value = "Result"
$ts.Range("C2").Value = "Result"
$ts.Range("D2").Value = "FAIL"
$ts.Range("E2").Value = "Call"

$ts.Range("A1:E2").Font.Bold = $true
$ts.Range("A1:E2").HorizontalAlignment = -4108

# Row 3: qlTimeSeries - currently failing (expected ts#0016, got ts#0000)
$ts.Range("A3").Value = "qlTimeSeries"
$ts.Range("B3").Value = "ts#0016"
$ts.Range("C3").Value = "ts#0000"
$ts.Range("D3").Formula = '=IF(B3=C3,"PASS","FAIL")'
$ts.Range("E3").Value = "ts#0000"

# Row 4: qlTimeSeriesFirstDate
$ts.Range("A4").Value = "qlTimeSeriesFirstDate"
$ts.Range("B4").Value = 25569
$ts.Range("C4").Value = 25569
$ts.Range("D4:D11").Formula = '=IF(B4=C4,"PASS","FAIL")'
$ts.Range("E4").Value = 25569

# Row 5: qlTimeSeriesLastDate
$ts.Range("A5").Value = "qlTimeSeriesLastDate"
$ts.Range("B5").Value = 25571
$ts.Range("C5").Value = 25571
$ts.Range("E5").Value = 25571

# Row 6: qlTimeSeriesSize
$ts.Range("A6").Value = "qlTimeSeriesSize"
$ts.Range("B6").Value = 3
$ts.Range("C6").Value = 3
$ts.Range("E6").Value = 3

# Row 7: qlTimeSeriesEmpty
$ts.Range("A7").Value = "qlTimeSeriesEmpty"
$ts.Range("B7").Value = $false
$ts.Range("C7").Value = $false
$ts.Range("E7").Value = $false

# Row 8: qlTimeSeriesDates
$ts.Range("A8").Value = "qlTimeSeriesDates"
$ts.Range("B8").Value = 25569
$ts.Range("C8").Value = 25569
$ts.Range("E8").Value = 25569

# Row 9: qlTimeSeriesValues
$ts.Range("A9").Value = "qlTimeSeriesValues"
$ts.Range("B9").Value = 1
$ts.Range("C9").Value = 1
$ts.Range("E9").Value = 1

# Row 10: qlTimeSeriesValue
$ts.Range("A10").Value = "qlTimeSeriesValue"
$ts.Range("B10").Value = 1
$ts.Range("C10").Value = 1
$ts.Range("E10").Value = 1

# Row 11: qlTimeSeriesFromIndex - currently erroring
$ts.Range("A11").Value = "qlTimeSeriesFromIndex"
$ts.Range("B11").Value = "#NUM!"
$ts.Range("C11").Value = "#NAME?"
$ts.Range("E11").Value = "#NAME?"

$ts.Range("D3:D11").HorizontalAlignment = -4108

$ts.Activate()
$ts.Range("A11:E11").Select()

# ---------------------------------------------------------------------
# 3. List the new group on the "UnitTests" summary sheet.
# ---------------------------------------------------------------------
$unitTests.Cells.Item(6, 1).Value = "Utilities"

# ---------------------------------------------------------------------
# 4. Restore the "Date" sheet's selection (unrelated to the new sheet,
#    but its scroll/selection shifts down two rows in the target file).
# ---------------------------------------------------------------------
$dateSheet.Activate()
$dateSheet.Range("A42:E42").Select()

# ---------------------------------------------------------------------
# 5. Utilities becomes the active tab again (now the last sheet).
# ---------------------------------------------------------------------
$utilSheet.Activate()
$utilSheet.Range("A3:E3").Select()
